$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a new, numeric-looking price that must stay plain text
# (matches the source data, which stores every price as a string).
$textCells = 'D5','D6','D7','D8','D10','D16','D18','D19','D20','D21','D22','D23','D25','D26','D28','D29','D33','D35','D36','D37','D38','D39','D40','D42','D45','D46','D47','D48','D49','D50'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '60.605.80'
$ws.Range("E2").Value = '  +3.60%  '

# Row 3
$ws.Range("D3").Value = '2.648.77'
$ws.Range("E3").Value = '  +1.10%  '

# Row 4
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").Value = '569.68'
$ws.Range("E5").Value = '  +6.44%  '

# Row 6
$ws.Range("D6").Value = '147.19'
$ws.Range("E6").Value = '  +3.07%  '

# Row 7
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.36%  '

# Row 8
$ws.Range("D8").Value = '0.611'
$ws.Range("E8").Value = '  +7.66%  '

# Row 9
$ws.Range("D9").Value = '2.675.02'
$ws.Range("E9").Value = '  +1.90%  '

# Row 10
$ws.Range("D10").Value = '6.84'
$ws.Range("E10").Value = '  -0.54%  '

# Row 11
$ws.Range("E11").Value = '  +4.59%  '

# Row 12
$ws.Range("E12").Value = '  +6.44%  '

# Row 13
$ws.Range("E13").Value = '  +2.87%  '

# Row 14
$ws.Range("D14").Value = '3.120.68'
$ws.Range("E14").Value = '  +1.25%  '

# Row 15
$ws.Range("D15").Value = '60.612.54'
$ws.Range("E15").Value = '  +3.76%  '

# Row 16
$ws.Range("D16").Value = '21.90'
$ws.Range("E16").Value = '  +5.88%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.698.79'
$ws.Range("E17").Value = '  +3.04%  '

# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.0000138'
$ws.Range("E18").Value = '  +4.79%  '

# Row 19
$ws.Range("D19").Value = '4.57'
$ws.Range("E19").Value = '  +3.75%  '

# Row 20
$ws.Range("D20").Value = '344.51'
$ws.Range("E20").Value = '  +3.22%  '

# Row 21
$ws.Range("D21").Value = '10.49'
$ws.Range("E21").Value = '  +3.69%  '

# Row 22
$ws.Range("D22").Value = '6.38'
$ws.Range("E22").Value = '  +2.59%  '

# Row 23
$ws.Range("D23").Value = '5.83'
$ws.Range("E23").Value = '  +1.36%  '

# Row 24
$ws.Range("E24").Value = '  -0.11%  '

# Row 25
$ws.Range("D25").Value = '66.67'
$ws.Range("E25").Value = '  +0.63%  '

# Row 26
$ws.Range("D26").Value = '0.442'
$ws.Range("E26").Value = '  +6.59%  '

# Row 27
$ws.Range("E27").Value = '  +1.91%  '

# Row 28
$ws.Range("D28").Value = '0.996'
$ws.Range("E28").Value = '  -0.22%  '

# Row 29
$ws.Range("D29").Value = '7.43'
$ws.Range("E29").Value = '  +4.83%  '

# Row 30
$ws.Range("E30").Value = '  +7.44%  '

# Row 31
$ws.Range("E31").Value = '  -0.03%  '

# Row 32
$ws.Range("E32").Value = '  +5.00%  '

# Row 33
$ws.Range("D33").Value = '6.21'
$ws.Range("E33").Value = '  +6.09%  '

# Row 34
$ws.Range("E34").Value = '  +2.25%  '

# Row 35
$ws.Range("D35").Value = '154.49'
$ws.Range("E35").Value = '  +2.50%  '

# Row 36
$ws.Range("D36").Value = '4.09'
$ws.Range("E36").Value = '  +5.53%  '

# Row 37
$ws.Range("D37").Value = '1.19'
$ws.Range("E37").Value = '  +8.34%  '

# Row 38
$ws.Range("D38").Value = '0.913'
$ws.Range("E38").Value = '  +12.43%  '

# Row 39
$ws.Range("D39").Value = '0.906'
$ws.Range("E39").Value = '  +7.00%  '

# Row 40
$ws.Range("D40").Value = '37.62'
$ws.Range("E40").Value = '  +1.41%  '

# Row 41
$ws.Range("E41").Value = '  +7.79%  '

# Row 42
$ws.Range("D42").Value = '305.85'
$ws.Range("E42").Value = '  +8.84%  '

# Row 43
$ws.Range("E43").Value = '  +3.13%  '

# Row 44
$ws.Range("E44").Value = '  -0.54%  '

# Row 45
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '0.0985'
$ws.Range("E45").Value = '  +5.48%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.609'
$ws.Range("E46").Value = '  +2.66%  '

# Row 47
$ws.Range("D47").Value = '0.0550'
$ws.Range("E47").Value = '  +4.63%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '129.15'
$ws.Range("E48").Value = '  +13.58%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '19.56'
$ws.Range("E49").Value = '  +3.32%  '

# Row 50
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '10.68'
$ws.Range("E50").Value = '  -0.12%  '

# Row 51
$ws.Range("E51").Value = '  +5.35%  '

